$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out erroneous placeholder values
$ws.Range("C3").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()

# Fix a typo'd name
$ws.Range("C6").Value = "ISAAC SANCHEZ"

# Fill in row 9 with parsed values instead of ERROR placeholders
$ws.Range("B9").Value = "Sunday_ February 19_ 2023"
$ws.Range("C9").Value = "SACKLYN"
$ws.Range("D9").Value = "4_00 pm"
$ws.Range("E9").Value = "10_00 pm"
